# The trading-sim re-run no longer prints each row to the terminal, but the
# recomputed profit/account_balance figures still need to land in the sheet.
# Column AF = profit, column AG = running account_balance (rows 2-45).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 45

$profitValues = @(
    216.78, -344.36, 216.78, 216.78, -150.68, 216.78, 216.78, -196.21, -88.17, 216.77, 216.76, 216.76, 216.76, 216.75, 108.38, 216.76, 108.38, -215.36, -109.49, -102.26, -136.59, -135.14, -265.95, 216.76, 216.77, -264.87, -231.62, 216.76, 216.76, 216.76, 108.38, -85.28, 216.76, 216.76, -129.36, 108.38, 216.76, 216.76, 216.76, 216.76, 108.38, 216.76, 108.38, -129.72
)

$balanceValues = @(
    816.78, 472.42, 689.1999999999999, 905.9799999999999, 755.3, 972.0799999999999, 1188.86, 992.6499999999999, 904.4799999999999, 1121.25, 1338.01, 1554.77, 1771.53, 1988.28, 2096.66, 2313.42, 2421.8, 2206.44, 2096.95, 1994.69, 1858.1, 1722.96, 1457.01, 1673.77, 1890.54, 1625.670000000001, 1394.050000000001, 1610.810000000001, 1827.570000000001, 2044.330000000001, 2152.71, 2067.43, 2284.190000000001, 2500.950000000001, 2371.590000000001, 2479.970000000001, 2696.73, 2913.490000000001, 3130.250000000001, 3347.010000000001, 3455.390000000001, 3672.150000000001, 3780.530000000002, 3650.810000000002
)

for ($i = 0; $i -lt $profitValues.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 32).Value = $profitValues[$i]   # column AF
    $ws.Cells.Item($row, 33).Value = $balanceValues[$i]  # column AG
}
